$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12: swap the "B" operand bit group (I12/J12) to match the "A"/"C" pattern used elsewhere.
$ws.Range("I12").Value = "000"
$ws.Range("J12").Value = "0"

# Row 15: swap the "A" operand bit group (D15/F15) and the "B" operand bit group (I15/K15).
# A leading apostrophe keeps these text cells on their original (quote-prefixed) cell style,
# matching the source file where only the <v> value changes, not the style index.
$ws.Range("D15").Value = "'011"
$ws.Range("F15").Value = "'1"
$ws.Range("I15").Value = "'000"
$ws.Range("K15").Value = "'0"

# Update the selection shown when the sheet was last saved.
$ws.Range("S2:S15").Select()
